$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Replace the data table (header row kept, data rows 2-6 updated with the
# new set of page / title records).
# ---------------------------------------------------------------------------
$data = @(
    @("path", "type", "language", "browser_title", "content_title"),
    @("/about-cancer/coping/feelings", "Article", "English", "Feelings and Cancer | CGDP - Dev", "Feelings and Cancer "),
    @("/espanol/cancer/sobrellevar/sentimientos/relajarse", "Article", "Spanish", "Aprenda a relajarse | CGDP - Dev", "Aprenda a relajarse"),
    @("/about-cancer/coping/feelings/relaxation", "Article", "English", "Learning to Relax | CGDP - Dev", "Learning to Relax"),
    @("/about-cancer/coping/feelings/relaxation/dfharvard", "Cancer Center", "English", "Dana Farber/Harvard Cancer Center | CGDP - Dev", "Dana Farber/Harvard Cancer Center"),
    @("/about-cancer/coping/feelings/relaxation/loukissas-jennifer", "Biography", "English", "Jennifer K. Loukissas, M.P.P. | CGDP - Dev", "Jennifer K. Loukissas, M.P.P.")
)

for ($r = 1; $r -le 6; $r++) {
    for ($c = 1; $c -le 5; $c++) {
        $ws.Cells.Item($r, $c).Value = $data[$r - 1][$c - 1]
    }
}

# ---------------------------------------------------------------------------
# Row 5 (Dana Farber/Harvard Cancer Center), column E (content_title) gets a
# new font color (dark gray, #222222) applied to it.
# ---------------------------------------------------------------------------
$ws.Range("E5").Font.Color = 2236962

# ---------------------------------------------------------------------------
# Column widths.
# ---------------------------------------------------------------------------
$ws.Columns.Item(1).ColumnWidth = 36.041666666666664
$ws.Columns.Item(2).ColumnWidth = 12.541666666666666
$ws.Columns.Item(4).ColumnWidth = 35.041666666666664
$ws.Columns.Item(5).ColumnWidth = 12.541666666666666

# ---------------------------------------------------------------------------
# Selection / active cell as last saved by the editing user.
# ---------------------------------------------------------------------------
$ws.Range("H12").Select() | Out-Null

# ---------------------------------------------------------------------------
# Page setup orientation (portrait).
# ---------------------------------------------------------------------------
$ws.PageSetup.Orientation = 1 | Out-Null
